$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footnote text clarified: the old note only mentioned "%", the new one
#     also covers the covariance-unit (%^2) case ---
$ws.Range("I10").Value = "Parametrin yksikkö on % tai kovarianssin yksikkö (%^2)"

# --- Data change: relative-view row 4 now references "DM Equities" instead
#     of "EM Equities" in the comparator (F4) ---
$ws.Range("F4").Value = "DM Equities"

# --- Formulas in column I (preview text) redefined:
#     * "corr" branch -> "cov" branch (the mean-rows only ever used the mean
#       branch, the corr branch was dead code left over from before cov views
#       existed)
#     * the "+1 offset" hack (E+1&"...") replaced by a clean string
#       concatenation (E&" + ...") now that weighted/relative views make
#       more sense
#     Setting the whole I2:I5 block in one shot mirrors how these rows share
#     one underlying formula pattern. ---
$ws.Range("I2:I5").Formula = "=IF(A2=""mean"",""Mean('""&B2&""') ""&D2&"" ""&IF(ISBLANK(F2),E2,E2&"" + Mean('""&F2&""')""),) & IF(A2=""cov"",""Cov('""&B2&""', '""&C2&""') ""&D2&"" ""&IF(ISBLANK(F2),E2,E2&"" + Cov('""&F2&""', '""&G2&""')""),)"

# --- Selection moved from G8 to C7 ---
$ws.Range("C7").Select() | Out-Null
